$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# Overview sheet: row 3 (ea62922e...md) - Status columns B3 (zh-cn) and C3 (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# zh-cn sheet: row 3 - Status (B3) and Latest Handback DateTime (G3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $status
$wsZhCn.Range("G3").Value = "2016-01-17 16:16:59"

# de-de sheet: row 3 - Status (B3) and Latest Handback DateTime (G3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $status
$wsDeDe.Range("G3").Value = "2016-01-17 16:17:19"
